$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.954.01"
$ws.Range("E2").Value = "  +5.65%  "
$ws.Range("D3").Value = "2.373.85"
$ws.Range("E3").Value = "  +4.08%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.55"
$ws.Range("E5").Value = "  +3.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.15"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "2.371.37"
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.337"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.19"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").Value = "2.793.78"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "60.842.98"
$ws.Range("E16").Value = "  +5.50%  "
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").Value = "2.381.57"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.79"
$ws.Range("E19").Value = "  +3.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.98"
$ws.Range("E20").Value = "  +10.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.21"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "318.32"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.62"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("E25").Value = "  +4.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.06"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.37"
$ws.Range("E28").Value = "  +6.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.76"
$ws.Range("E29").Value = "  +3.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.11"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "0.0₃0739"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +11.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.92"
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.43"
$ws.Range("E34").Value = "  +16.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.386"
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.14"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.22"
$ws.Range("E39").Value = "  +9.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "321.56"
$ws.Range("E40").Value = "  +12.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.55"
$ws.Range("E41").Value = "  +4.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.33"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "144.61"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.49"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0958"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.45"
$ws.Range("E46").Value = "  +8.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0502"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.566"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0215"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0211"
$ws.Range("E50").Value = "  +5.45%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.04"
$ws.Range("E51").Value = "  +1.02%  "
